# Generate Report for Handoff
# Adds a new row (row 8) describing the file
# f30d04b9-72a7-4a23-a900-886aec82188b.md to the Overview, zh-cn and de-de
# worksheets of the localization-status report, matching the newly
# "Ready for handoff" status of that file.

$wb = $excel.ActiveWorkbook

$fileName   = "f30d04b9-72a7-4a23-a900-886aec82188b.md"
$status     = "Ready for handoff"
$overviewDate = "2016-20-17 06:20:25"

$zhXlf      = "f30d04b9-72a7-4a23-a900-886aec82188b.1d79bc8fae4e2cf21ca9860ed7e97361c8642351.zh-cn.xlf"
$zhDate     = "2016-03-17 06:20:18"

$deXlf      = "f30d04b9-72a7-4a23-a900-886aec82188b.1d79bc8fae4e2cf21ca9860ed7e97361c8642351.de-de.xlf"
$deDate     = "2016-03-17 06:20:25"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/f2e18d56ac214308e2b2030aa92eb037813abc30/e2e/$fileName"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85011af4eedbe76c4ab1708a033b682a723c90c2/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$zhXlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5045a4f53f6337978e194d7e359e2b801bb708f1/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$deXlf"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$row = 8

$wsOverview.Hyperlinks.Add($wsOverview.Range("A$row"), $mdUrl, "", "", $fileName)
$wsOverview.Range("B$row").Value2 = $status
$wsOverview.Range("C$row").Value2 = $status
$wsOverview.Range("D$row").Value2 = $overviewDate

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A$row"), $mdUrl, "", "", $fileName)
$wsZh.Hyperlinks.Add($wsZh.Range("B$row"), $mdUrl, "", "", ".md")
$wsZh.Range("C$row").Value2 = $status
$wsZh.Hyperlinks.Add($wsZh.Range("D$row"), $zhXlfUrl, "", "", $zhXlf)
$wsZh.Range("E$row").Value2 = $zhDate
$wsZh.Range("H$row").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("I$row").Value2 = "Include"

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A$row"), $mdUrl, "", "", $fileName)
$wsDe.Hyperlinks.Add($wsDe.Range("B$row"), $mdUrl, "", "", ".md")
$wsDe.Range("C$row").Value2 = $status
$wsDe.Hyperlinks.Add($wsDe.Range("D$row"), $deXlfUrl, "", "", $deXlf)
$wsDe.Range("E$row").Value2 = $deDate
$wsDe.Range("H$row").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("I$row").Value2 = "Include"
